$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.367.99"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "3.689.86"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "680.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.46"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  -3.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.440"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").Value = "4.313.31"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "3.690.16"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "69.348.13"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.05"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.94"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "3.836.11"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.14"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.69"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("E30").Value = "  -4.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.97"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").Value = "3.677.93"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("E36").Value = "  -6.91%  "
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.26"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E40").Value = "  -4.56%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0907"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "170.73"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.944"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.61"
$ws.Range("D45").ClearFormats()
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.45"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.47%  "
$ws.Range("B47").Value = "SuiNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.70"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.80"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.04%  "
